$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with the latest crypto snapshot.
# D-column values are written via a literal-text formula + copy/paste-values so
# Excel keeps them as plain text (preserving things like trailing zeros, the
# "12.34.56"-style price strings, and the subscript-6 notation) instead of
# auto-converting them to numbers.

$ws.Range("D2").Formula = "=`"68.115.18`""
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("E2").Value = '  -1.06%  '
$ws.Range("D3").Formula = "=`"2.641.86`""
$ws.Range("D3").Copy()
$ws.Range("D3").PasteSpecial(-4163)
$ws.Range("E3").Value = '  -0.31%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Formula = "=`"596.75`""
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = '  -0.70%  '
$ws.Range("D6").Formula = "=`"155.46`""
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = '  -0.34%  '
$ws.Range("E8").Value = '  -1.13%  '
$ws.Range("D9").Formula = "=`"0.140`""
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = '  -0.37%  '
$ws.Range("E10").Value = '  -0.99%  '
$ws.Range("E11").Value = '  -0.16%  '
$ws.Range("E12").Value = '  -0.43%  '
$ws.Range("D13").Formula = "=`"27.91`""
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = '  -0.28%  '
$ws.Range("E14").Value = '  -0.53%  '
$ws.Range("D15").Formula = "=`"3.123.75`""
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = '  -0.27%  '
$ws.Range("D16").Formula = "=`"68.074.50`""
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = '  -0.92%  '
$ws.Range("D17").Formula = "=`"2.661.06`""
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Value = '  +0.53%  '
$ws.Range("D18").Formula = "=`"11.33`""
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = '  -0.50%  '
$ws.Range("D19").Formula = "=`"362.89`""
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E20").Value = '  -1.26%  '
$ws.Range("E21").Value = '  +2.61%  '
$ws.Range("D22").Formula = "=`"4.77`""
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = '  -3.06%  '
$ws.Range("E23").Value = '  -2.54%  '
$ws.Range("D24").Formula = "=`"74.59`""
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = '  +2.37%  '
$ws.Range("E25").Value = '  -0.04%  '
$ws.Range("D26").Formula = "=`"9.71`""
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = '  -3.81%  '
$ws.Range("D27").Formula = "=`"2.773.79`""
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = '  -0.13%  '
$ws.Range("E28").Value = '  -2.47%  '
$ws.Range("D29").Formula = "=`"0.999`""
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = '  -0.12%  '
$ws.Range("D30").Formula = "=`"554.82`""
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = '  -5.05%  '
$ws.Range("E31").Value = '  -0.33%  '
$ws.Range("E32").Value = '  -1.71%  '
$ws.Range("E33").Value = '  -1.08%  '
$ws.Range("E34").Value = '  -2.66%  '
$ws.Range("D35").Formula = "=`"0.999`""
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("E36").Value = '  -0.91%  '
$ws.Range("D37").Formula = "=`"161.18`""
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = '  +0.51%  '
$ws.Range("D38").Formula = "=`"19.41`""
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = '  +0.28%  '
$ws.Range("E39").Value = '  +0.84%  '
$ws.Range("E40").Value = '  -4.02%  '
$ws.Range("E41").Value = '  -1.51%  '
$ws.Range("D42").Formula = "=`"0.0₆0335`""
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = '  +4.18%  '
$ws.Range("D43").Formula = "=`"17.80`""
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = '  +0.40%  '
$ws.Range("D44").Formula = "=`"2.60`""
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = '  -2.52%  '
$ws.Range("E45").Value = '  +0.02%  '
$ws.Range("D46").Formula = "=`"159.44`""
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("D47").Formula = "=`"3.71`""
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = '  -0.69%  '
$ws.Range("D48").Formula = "=`"22.00`""
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = '  -0.46%  '
$ws.Range("E49").Value = '  -1.53%  '
$ws.Range("E50").Value = '  -0.25%  '
$ws.Range("E51").Value = '  -0.70%  '

$excel.CutCopyMode = $false
